$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking text values to stay as text (matches their existing
# General-but-string storage) by setting NumberFormat to Text before assigning.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cryptos list values (Price / Volume(1h) columns).
$ws.Range("D2").Value = "66.882.05"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.282.35"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "571.84"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "174.50"
$ws.Range("E6").Value = "  -5.73%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "3.277.00"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "45.32"
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "685.14"
$ws.Range("E14").Value = "  +6.61%  "
$ws.Range("D15").Value = "3.809.88"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "8.26"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "66.983.12"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "3.281.72"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").Value = "17.25"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").Value = "0.885"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "17.02"
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").Value = "5.17"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("D25").Value = "99.57"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").Value = "33.55"
$ws.Range("E29").Value = "  +7.53%  "
$ws.Range("D30").Value = "8.37"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "573.49"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").Value = "3.871.16"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "55.33"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -13.99%  "
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").Value = "0.0₃0667"
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  -5.77%  "
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +6.42%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "129.80"
$ws.Range("E51").Value = "  -0.62%  "
